$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.844.69"
$ws.Range("E2").Value = "'  -0.31%  "
$ws.Range("D3").Value = "'1.585.21"
$ws.Range("E3").Value = "'  -2.10%  "
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("D5").Value = "'209.40"
$ws.Range("E5").Value = "'  -1.63%  "
$ws.Range("E6").Value = "'  +0.00%  "
$ws.Range("D7").Value = "'0.482"
$ws.Range("E7").Value = "'  -3.48%  "
$ws.Range("E8").Value = "'  -0.87%  "
$ws.Range("D9").Value = "'0.0616"
$ws.Range("E9").Value = "'  -0.28%  "
$ws.Range("D10").Value = "'18.07"
$ws.Range("E10").Value = "'  -1.66%  "
$ws.Range("D11").Value = "'0.0790"
$ws.Range("D12").Value = "'1.806.08"
$ws.Range("E12").Value = "'  -2.04%  "
$ws.Range("D13").Value = "'1.579.89"
$ws.Range("E13").Value = "'  -2.49%  "
$ws.Range("D14").Value = "'4.03"
$ws.Range("E14").Value = "'  -2.49%  "
$ws.Range("D15").Value = "'0.510"
$ws.Range("E15").Value = "'  -2.65%  "
$ws.Range("D16").Value = "'25.832.95"
$ws.Range("E16").Value = "'  -0.37%  "
$ws.Range("D17").Value = "'0.0₃0723"
$ws.Range("E17").Value = "'  -1.76%  "
$ws.Range("D18").Value = "'59.86"
$ws.Range("E18").Value = "'  -2.77%  "
$ws.Range("E19").Value = "'  +0.05%  "
$ws.Range("D20").Value = "'191.09"
$ws.Range("E20").Value = "'  -0.23%  "
$ws.Range("D21").Value = "'4.19"
$ws.Range("E21").Value = "'  -1.41%  "
$ws.Range("D22").Value = "'9.37"
$ws.Range("E22").Value = "'  -1.41%  "
$ws.Range("D23").Value = "'5.92"
$ws.Range("E23").Value = "'  -1.95%  "
$ws.Range("E24").Value = "'  -1.02%  "
$ws.Range("D25").Value = "'141.54"
$ws.Range("E25").Value = "'  -1.59%  "
$ws.Range("E26").Value = "'  -0.04%  "
$ws.Range("E27").Value = "'  -1.22%  "
$ws.Range("D28").Value = "'15.09"
$ws.Range("E28").Value = "'  -0.94%  "
$ws.Range("D29").Value = "'6.44"
$ws.Range("E29").Value = "'  -2.96%  "
$ws.Range("E30").Value = "'  -5.61%  "
$ws.Range("E31").Value = "'  -1.05%  "
$ws.Range("D32").Value = "'3.11"
$ws.Range("E32").Value = "'  -0.47%  "
$ws.Range("E33").Value = "'  -2.39%  "
$ws.Range("E34").Value = "'  +0.18%  "
$ws.Range("E35").Value = "'  -2.42%  "
$ws.Range("D36").Value = "'1.098.29"
$ws.Range("E36").Value = "'  -2.38%  "
$ws.Range("E37").Value = "'  +0.07%  "
$ws.Range("D38").Value = "'2.33"
$ws.Range("E38").Value = "'  -2.20%  "
$ws.Range("D39").Value = "'0.506"
$ws.Range("E39").Value = "'  -1.22%  "
$ws.Range("E40").Value = "'  -2.03%  "
$ws.Range("B41").Value = "'TrustWalletToken"
$ws.Range("C41").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.821"
$ws.Range("E41").Value = "'  +9.98%  "
$ws.Range("B42").Value = "'ARBITRUM"
$ws.Range("C42").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'0.779"
$ws.Range("E42").Value = "'  -7.48%  "
$ws.Range("D43").Value = "'5.20"
$ws.Range("E43").Value = "'  +2.60%  "
$ws.Range("D44").Value = "'93.78"
$ws.Range("E44").Value = "'  -4.11%  "
$ws.Range("D45").Value = "'1.718.92"
$ws.Range("E45").Value = "'  -2.05%  "
$ws.Range("E46").Value = "'  +0.16%  "
$ws.Range("E47").Value = "'  -0.94%  "
$ws.Range("D48").Value = "'53.19"
$ws.Range("E48").Value = "'  -1.61%  "
$ws.Range("D49").Value = "'0.0509"
$ws.Range("E49").Value = "'  -1.63%  "
$ws.Range("E50").Value = "'  -0.80%  "
$ws.Range("E51").Value = "'  -0.08%  "
